# Apply updated "Inscritos"/"Pagos"/"Inscrições homologadas" counts
# to the Inscricoes worksheet, as reflected in the source diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

$ws.Range("E4").Value  = 24

$ws.Range("E8").Value  = 15

$ws.Range("E12").Value = 3

$ws.Range("E15").Value = 160
$ws.Range("F15").Value = 85
$ws.Range("H15").Value = 126

$ws.Range("E17").Value = 121
$ws.Range("F17").Value = 58
$ws.Range("H17").Value = 90

$ws.Range("E18").Value = 110

$ws.Range("E24").Value = 23
$ws.Range("F24").Value = 14
$ws.Range("H24").Value = 18

$ws.Range("E26").Value = 28
$ws.Range("F26").Value = 14
$ws.Range("H26").Value = 24

$ws.Range("E33").Value = 42
$ws.Range("F33").Value = 13
$ws.Range("H33").Value = 25

$ws.Range("E36").Value = 100

$ws.Range("E37").Value = 54

$ws.Range("F41").Value = 16
$ws.Range("H41").Value = 27

$ws.Range("E49").Value = 72

$ws.Range("E51").Value = 11

$ws.Range("E52").Value = 7

$ws.Range("E57").Value = 14

$ws.Range("E67").Value = 39

$ws.Range("E71").Value = 36

$ws.Range("E74").Value = 21

$ws.Range("E75").Value = 15

$ws.Range("E78").Value = 46
$ws.Range("F78").Value = 21
$ws.Range("H78").Value = 42

$ws.Range("E88").Value = 23
